$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("C5").Value = 52

# Row 6
$ws.Range("B6").Value = "<people>"
$ws.Range("C6").Value = 58

# Row 7
$ws.Range("C7").Value = 52

# Row 8
$ws.Range("B8").Value = "<was>"
$ws.Range("C8").Value = 55

# Row 9
$ws.Range("B9").Value = "<word>"
$ws.Range("C9").Value = 48

# Row 10
$ws.Range("B10").Value = "<be>"
$ws.Range("C10").Value = 58

# Row 11
$ws.Range("B11").Value = "<echo>"
$ws.Range("C11").Value = 54

# Row 12
$ws.Range("B12").Value = "<yes>"
$ws.Range("C12").Value = 52

# Row 13
$ws.Range("B13").Value = "<omward>"
$ws.Range("C13").Value = 52

# Row 14
$ws.Range("C14").Value = 56

# Row 15
$ws.Range("C15").Value = 53

# Row 16
$ws.Range("B16").Value = "<their>"
$ws.Range("C16").Value = 30
